$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy header style from an existing header cell (e.g. A1) so AC1:AE1 match formatting
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

# Fill data rows 2-46 with team record values
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 29).Value = 65
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 0
}
